# becker menus.xlsx -- fix sheet (tab) names so every sheet uses a
# consistent, unambiguous MMDDYYYY date format.
#
# Several sheet names were missing the leading zero on single-digit
# months/days (e.g. "1112020" meant to be 11/01/2020, "542021" meant to be
# 05/04/2021). The sheet "10142021" was created by mistake and was left
# completely empty -- it was really meant to be "12042021", and every sheet
# after it in that block was off by one slot for the same reason (ending in
# a genuine duplicate of "12102021" at the very end). Removing that one
# empty/misnamed sheet and renaming the rest of the block forward realigns
# everything without touching any of the real menu data.

$wb = $excel.ActiveWorkbook

# The empty, mis-dated placeholder sheet -- its data never got entered
# because the whole block that followed it was shifted by one name.
[void]$wb.Worksheets.Item("10142021").Delete()

# 2020 sheets missing the leading zero for November.
$wb.Worksheets.Item("1112020").Name = "11012020"
$wb.Worksheets.Item("1122020").Name = "11022020"
$wb.Worksheets.Item("1132020").Name = "11032020"
$wb.Worksheets.Item("1142020").Name = "11042020"

# 2021 sheets missing the leading zero for April.
$wb.Worksheets.Item("4132021").Name = "04132021"
$wb.Worksheets.Item("4142021").Name = "04142021"
$wb.Worksheets.Item("4152021").Name = "04152021"
$wb.Worksheets.Item("4162021").Name = "04162021"
$wb.Worksheets.Item("4172021").Name = "04172021"
$wb.Worksheets.Item("4182021").Name = "04182021"
$wb.Worksheets.Item("4192021").Name = "04192021"
$wb.Worksheets.Item("4202021").Name = "04202021"
$wb.Worksheets.Item("4212021").Name = "04212021"
$wb.Worksheets.Item("4222021").Name = "04222021"

# 2021 sheets missing the leading zero for May (and, for the first six, the
# leading zero on the day too).
$wb.Worksheets.Item("542021").Name = "05042021"
$wb.Worksheets.Item("552021").Name = "05052021"
$wb.Worksheets.Item("562021").Name = "05062021"
$wb.Worksheets.Item("572021").Name = "05072021"
$wb.Worksheets.Item("582021").Name = "05082021"
$wb.Worksheets.Item("592021").Name = "05092021"
$wb.Worksheets.Item("5102021").Name = "05102021"
$wb.Worksheets.Item("5112021").Name = "05112021"
$wb.Worksheets.Item("5122021").Name = "05122021"
$wb.Worksheets.Item("5242021").Name = "05242021"
$wb.Worksheets.Item("5252021").Name = "05252021"
$wb.Worksheets.Item("5262021").Name = "05262021"
$wb.Worksheets.Item("5272021").Name = "05272021"
$wb.Worksheets.Item("5282021").Name = "05282021"
$wb.Worksheets.Item("5292021").Name = "05292021"
$wb.Worksheets.Item("5302021").Name = "05302021"

# 2021 sheet missing the leading zero for September.
$wb.Worksheets.Item("9172021").Name = "09172021"

# 2021 sheets missing the leading zero for October.
$wb.Worksheets.Item("1072021").Name = "10072021"
$wb.Worksheets.Item("1082021").Name = "10082021"
$wb.Worksheets.Item("1092021").Name = "10092021"

# The rest of the December block shifts forward one slot to fill the gap
# left by the empty "10142021" sheet above, and picks up its proper
# leading-zero month/day along the way.
$wb.Worksheets.Item("1242021").Name = "12042021"
$wb.Worksheets.Item("1252021").Name = "12052021"
$wb.Worksheets.Item("1262021").Name = "12062021"
$wb.Worksheets.Item("1272021").Name = "12072021"
$wb.Worksheets.Item("1282021").Name = "12082021"
$wb.Worksheets.Item("1292021").Name = "12092021"

# The old final sheet "12102021" was an exact duplicate of the date that
# "1292021" now correctly becomes, so it keeps its name unchanged and no
# further action is needed for it.
